# Apply the updated cryptocurrency price/volume snapshot to Sheet1.
# Rows 24/25, 26/27 and 47/48 additionally swap rank order (name + link
# move together with the row's price/volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.983.55'
$ws.Range("E2").Value = '  +2.87%  '
$ws.Range("D3").Value = '2.968.98'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.96'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.02'
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").Value = '2.967.00'
$ws.Range("E8").Value = '  +1.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.507'
$ws.Range("E9").Value = '  +0.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.24'
$ws.Range("E10").Value = '  +3.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.152'
$ws.Range("E11").Value = '  +6.27%  '
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("E13").Value = '  +7.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.41'
$ws.Range("E14").Value = '  -1.03%  '
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '3.463.24'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").Value = '62.852.02'
$ws.Range("E17").Value = '  +2.58%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.76'
$ws.Range("E18").Value = '  +0.31%  '
$ws.Range("D19").Value = '2.972.93'
$ws.Range("E19").Value = '  +1.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '443.01'
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.56'
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.674'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.11'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("B24").Value = 'RenderToken'
$ws.Range("C24").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.33'
$ws.Range("E24").Value = '  +2.22%  '
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.70'
$ws.Range("E25").Value = '  -0.34%  '
$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.15'
$ws.Range("E26").Value = '  -3.09%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.92'
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.28'
$ws.Range("E29").Value = '  +4.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.62'
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.17'
$ws.Range("E31").Value = '  -3.30%  '
$ws.Range("D32").Value = '0.0₃0980'
$ws.Range("E32").Value = '  +11.44%  '
$ws.Range("E33").Value = '  -0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '26.66'
$ws.Range("E34").Value = '  -0.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("E36").Value = '  -1.65%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.68'
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.11'
$ws.Range("E38").Value = '  +4.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.05'
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.54'
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.56'
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.119'
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.283'
$ws.Range("E43").Value = '  -0.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.87'
$ws.Range("E44").Value = '  -1.26%  '
$ws.Range("D45").Value = '2.745.67'
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '134.82'
$ws.Range("E46").Value = '  +1.00%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '367.32'
$ws.Range("E47").Value = '  -2.57%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0341'
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.16'
$ws.Range("E50").Value = '  -3.34%  '
$ws.Range("E51").Value = '  -0.37%  '
